$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-15
# from serial date 45204 (2023-10-05) to 45207 (2023-10-08),
# preserving the existing date number format.
$ws.Range("C2:C15").Value = 45207
